$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; existing rows 7..62 shift down to 8..63.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new record's data.
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 44552
$ws.Cells.Item(7, 4).Style = $ws.Cells.Item(8, 4).Style
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = 100112021
$ws.Cells.Item(7, 7).Value = "Ají"
$ws.Cells.Item(7, 8).Value = "Americana (o)"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 40
$ws.Cells.Item(7, 11).Value = 36000
$ws.Cells.Item(7, 12).Value = 38000
$ws.Cells.Item(7, 13).Value = 37000
$ws.Cells.Item(7, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 1480
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"
